# Update "想去人数" (interested-count) values across sheets to reflect
# newly generated stats (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 415
$ws1.Range("F3").Value = 2495
$ws1.Range("F4").Value = 116

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 114
$ws2.Range("F3").Value = 2

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 415
$ws4.Range("F3").Value = 114
$ws4.Range("F4").Value = 2
$ws4.Range("F7").Value = 2495
$ws4.Range("F8").Value = 116
